$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 66
$ws1.Range("F3").Value = 11691
$ws1.Range("F4").Value = 214
$ws1.Range("F5").Value = 337
$ws1.Range("F7").Value = 11643
$ws1.Range("F8").Value = 481
$ws1.Range("F10").Value = 94
$ws1.Range("F11").Value = 1765
$ws1.Range("F12").Value = 5764
$ws1.Range("F14").Value = 3513
$ws1.Range("F16").Value = 17

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 570

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 570
$ws4.Range("F3").Value = 66
$ws4.Range("F5").Value = 11691
$ws4.Range("F6").Value = 214
$ws4.Range("F7").Value = 337
$ws4.Range("F9").Value = 11643
$ws4.Range("F10").Value = 481
$ws4.Range("F12").Value = 94
$ws4.Range("F13").Value = 1765
$ws4.Range("F15").Value = 5764
$ws4.Range("F17").Value = 3513
$ws4.Range("F19").Value = 17
